$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $txt = $shp.TextFrame.TextRange.Text
        if ($txt -eq "Decentralized") {
            $shp.TextFrame.TextRange.Text = "Shared"
        } elseif ($txt -eq "Centralized") {
            $shp.TextFrame.TextRange.Text = "Singular"
        }
    }
}
